$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 137 (shifts existing rows 137-144 down to 138-145)
$ws.Rows.Item(137).Insert()

# Match the style of the date cell (D column uses a date-formatted style) to the
# rest of the data block, without touching the whole row's formatting.
$ws.Range("D137").NumberFormat = $ws.Range("D138").NumberFormat

# Populate the new row 137 with the new weekly record
$ws.Range("A137").Value = 11
$ws.Range("B137").Value = "Vega Monumental Concepción"
$ws.Range("C137").Value = "Bíobío"
$ws.Range("D137").Value = 44610
$ws.Range("E137").Value = 8
$ws.Range("F137").Value = "Fruta"
$ws.Range("G137").Value = 100108
$ws.Range("H137").Value = "Tropicales y subtropicales"
$ws.Range("I137").Value = 100108005
$ws.Range("J137").Value = "Piña"
$ws.Range("K137").Value = "Caramelo"
$ws.Range("L137").Value = "Segunda"
$ws.Range("M137").Value = 150
$ws.Range("N137").Value = 15000
$ws.Range("O137").Value = 16000
$ws.Range("P137").Value = 15533
$ws.Range("Q137").Value = "$/caja 14 unidades"
$ws.Range("R137").Value = "Ecuador"
$ws.Range("S137").Value = 1110
$ws.Range("T137").Value = 14
